$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 99, shifting existing rows 99-112 down to 100-113.
$ws.Rows.Item(99).Insert()

# Populate the newly inserted row 99 with the new weekly record.
$ws.Cells.Item(99, 1).Value = 2
$ws.Cells.Item(99, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(99, 3).Value = "Coquimbo"
$ws.Cells.Item(99, 4).Value = 45218
$ws.Cells.Item(99, 5).Value = 4
$ws.Cells.Item(99, 6).Value = 100112022
$ws.Cells.Item(99, 7).Value = "Arveja Verde"
$ws.Cells.Item(99, 8).Value = "Perfection"
$ws.Cells.Item(99, 9).Value = "Primera"
$ws.Cells.Item(99, 10).Value = 400
$ws.Cells.Item(99, 11).Value = 20000
$ws.Cells.Item(99, 12).Value = 22000
$ws.Cells.Item(99, 13).Value = 21000
$ws.Cells.Item(99, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(99, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(99, 16).Value = 840
$ws.Cells.Item(99, 17).Value = 25
$ws.Cells.Item(99, 18).Value = "Hortaliza"
